# "WORKING: clean-up of code"
#
# 1. Rename the shared "OS-Drive" / "Data-Drive" labels (used as column
#    headers on several sheets) to "my OS-drive" / "my data-drive".
# 2. Add the (previously missing) header row to the "icons" sheet so it
#    also shows those two labels in A1/B1, formatted like the headers on
#    the other sheets.
# 3. Reset/update the remembered cell-selection on each sheet.

$wb = $excel.ActiveWorkbook

# --- 1. Rename the header labels everywhere they are used -----------------
foreach ($sheet in $wb.Worksheets) {
    $sheet.Cells.Replace("OS-Drive", "my OS-drive", 1, 1, $false, $false, $false, $false) | Out-Null
    $sheet.Cells.Replace("Data-Drive", "my data-drive", 1, 1, $false, $false, $false, $false) | Out-Null
}

# --- 2. Add the missing header row on the "icons" sheet --------------------
$wsFieldnames = $wb.Worksheets.Item("fieldnames")
$wsIcons = $wb.Worksheets.Item("icons")

$wsIcons.Range("A1").Value = "my OS-drive"
$wsIcons.Range("B1").Value = "my data-drive"

# Copy the header formatting from the "fieldnames" sheet onto the new cells.
$wsFieldnames.Range("A1:B1").Copy()
$wsIcons.Range("A1:B1").PasteSpecial(-4122)
$excel.CutCopyMode = $false

# --- 3. Update the remembered selection on each sheet -----------------------
$wsFieldnames.Activate()
$wsFieldnames.Range("A1").Select() | Out-Null

$wsUrl = $wb.Worksheets.Item("URL")
$wsUrl.Activate()
$wsUrl.Range("A1").Select() | Out-Null

$wsComments = $wb.Worksheets.Item("comments")
$wsComments.Activate()
$wsComments.Range("A1").Select() | Out-Null

$wsIcons.Activate()
$wsIcons.Range("B21").Select() | Out-Null
